# Claudminity cognitive-test-results update:
#  - consolidate old rows (drop the two empty-data timestamp rows)
#  - append newly collected trial rows (now 10 total data rows)
#  - add a "condition" column (H) marking each trial sober / THC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare date-formatted cells for the new rows by copying the
#     existing date style (s="1") from A4 down through A13, then
#     overwrite with the real values below. ---
$ws.Range("A4").Copy($ws.Range("A9:A13"))

# --- Row 4 ---
$ws.Range("A4").Value = 45464.983761574076
$ws.Range("B4").Value = 21.03
$ws.Range("C4").Value = 56.29
$ws.Range("D4").Value = 5.2
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = 10.38
$ws.Range("G4").Value = 51.89

# --- Row 5 ---
$ws.Range("A5").Value = 45465.918240740742
$ws.Range("B5").Value = 8.0399999999999991
$ws.Range("C5").Value = 40.18
$ws.Range("D5").Value = 4.59
$ws.Range("E5").Value = 22.94
$ws.Range("F5").Value = 7.99
$ws.Range("G5").Value = 39.97

# --- Row 6 ---
$ws.Range("A6").Value = 45472.964467592596
$ws.Range("B6").Value = 9.35
$ws.Range("C6").Value = 46.76
$ws.Range("D6").Value = 26.64
$ws.Range("E6").Value = 20.83
$ws.Range("F6").Value = 6.7
$ws.Range("G6").Value = 33.520000000000003

# --- Row 7 ---
$ws.Range("A7").Value = 45475.812881944446
$ws.Range("B7").Value = 7.29
$ws.Range("C7").Value = 36.450000000000003
$ws.Range("D7").Value = 6.36
$ws.Range("E7").Value = 31.79
$ws.Range("F7").Value = 8.69
$ws.Range("G7").Value = 43.44
$ws.Range("H7").Value = "sober"

# --- Row 8 ---
$ws.Range("A8").Value = 45475.978101851855
$ws.Range("B8").Value = 7.85
$ws.Range("C8").Value = 39.25
$ws.Range("D8").Value = 4.95
$ws.Range("E8").Value = 24.77
$ws.Range("F8").Value = 6.31
$ws.Range("G8").Value = 31.57
$ws.Range("H8").Value = "THC"

# --- Row 9 ---
$ws.Range("A9").Value = 45475.984467592592
$ws.Range("B9").Value = 8.57
$ws.Range("C9").Value = 42.83
$ws.Range("D9").Value = 5.32
$ws.Range("E9").Value = 26.61
$ws.Range("F9").Value = 7.76
$ws.Range("G9").Value = 38.82
$ws.Range("H9").Value = "THC"

# --- Row 10 ---
$ws.Range("A10").Value = 45475.989675925928
$ws.Range("B10").Value = 10.54
$ws.Range("C10").Value = 52.71
$ws.Range("D10").Value = 6.81
$ws.Range("E10").Value = 34.03
$ws.Range("F10").Value = 6.94
$ws.Range("G10").Value = 34.700000000000003
$ws.Range("H10").Value = "THC"

# --- Row 11 ---
$ws.Range("A11").Value = 45476.879537037035
$ws.Range("B11").Value = 7.37
$ws.Range("C11").Value = 36.86
$ws.Range("D11").Value = 4.0599999999999996
$ws.Range("E11").Value = 20.3
$ws.Range("F11").Value = 5.48
$ws.Range("G11").Value = 27.39
$ws.Range("H11").Value = "sober"

# --- Row 12 ---
$ws.Range("A12").Value = 45477.449652777781
$ws.Range("B12").Value = 5.27
$ws.Range("C12").Value = 26.35
$ws.Range("D12").Value = 4.84
$ws.Range("E12").Value = 24.19
$ws.Range("F12").Value = 5.17
$ws.Range("G12").Value = 25.84
$ws.Range("H12").Value = "sober"

# --- Row 13 ---
$ws.Range("A13").Value = 45477.525057870371
$ws.Range("B13").Value = 6.18
$ws.Range("C13").Value = 30.88
$ws.Range("D13").Value = 15.57
$ws.Range("E13").Value = 20.56
$ws.Range("F13").Value = 5.13
$ws.Range("G13").Value = 25.63
$ws.Range("H13").Value = "sober"

# --- Selection matches what was left selected after the edit ---
$ws.Range("A3:H13").Select() | Out-Null
